$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1031
$ws.Range("J17").Value = 1031
$ws.Range("L17").Value = 3093
$ws.Range("N17").Value = -3429

$ws.Range("H33").Value = 263.8625
$ws.Range("I33").Value = 199.98648
$ws.Range("K33").Value = 199.98648
$ws.Range("M33").Value = 29.01352

$ws.Range("H64").Value = 2868.75
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 2860
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 2860
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -3356

$ws.Range("H67").Value = 2868.75
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 2860
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 2860
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -4576

$ws.Range("H138").Value = 1428.1616
$ws.Range("I138").Value = 644.5417
$ws.Range("J138").Value = 2165.6863
$ws.Range("K138").Value = 1933.6251
$ws.Range("L138").Value = 6497.0589
$ws.Range("M138").Value = 3206.3749
$ws.Range("N138").Value = -16777.0589

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1427.4062
$ws.Range("I2").Value = 1564
$ws.Range("J2").Value = 835.5
$ws.Range("K2").Value = 1564
$ws.Range("L2").Value = 835.5
$ws.Range("M2").Value = -1451
$ws.Range("N2").Value = -1061.5

$ws.Range("H32").Value = 9645.402
$ws.Range("I32").Value = 8946.692999999999
$ws.Range("K32").Value = 8946.692999999999
$ws.Range("M32").Value = -8659.692999999999

$ws.Range("H74").Value = 1713.2273
$ws.Range("I74").Value = 1568.8889
$ws.Range("J74").Value = 2362.75
$ws.Range("K74").Value = 1568.8889
$ws.Range("L74").Value = 2362.75
$ws.Range("M74").Value = -694.8888999999999
$ws.Range("N74").Value = -4110.75

$ws.Range("H77").Value = 1713.2273
$ws.Range("I77").Value = 1568.8889
$ws.Range("J77").Value = 2362.75
$ws.Range("K77").Value = 7844.4445
$ws.Range("L77").Value = 11813.75
$ws.Range("M77").Value = -3476.4445
$ws.Range("N77").Value = -20549.75

$ws.Range("H102").Value = 14238.389
$ws.Range("I102").Value = 2480
$ws.Range("J102").Value = 32715.857
$ws.Range("K102").Value = 2480
$ws.Range("L102").Value = 32715.857
$ws.Range("M102").Value = -858
$ws.Range("N102").Value = -35959.857

$ws.Range("H110").Value = 1494.9354
$ws.Range("I110").Value = 1541
$ws.Range("J110").Value = 1362.5
$ws.Range("K110").Value = 1541
$ws.Range("L110").Value = 1362.5
$ws.Range("M110").Value = 504
$ws.Range("N110").Value = -5452.5

$ws.Range("H116").Value = 1427.4062
$ws.Range("I116").Value = 1564
$ws.Range("J116").Value = 835.5
$ws.Range("K116").Value = 1564
$ws.Range("L116").Value = 835.5
$ws.Range("M116").Value = 730
$ws.Range("N116").Value = -5423.5

$ws.Range("H132").Value = 12196861
$ws.Range("I132").Value = 15626085
$ws.Range("J132").Value = 4063.7778
$ws.Range("K132").Value = 46878255
$ws.Range("L132").Value = 12191.3334
$ws.Range("M132").Value = -46875725
$ws.Range("N132").Value = -17251.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1427.4062
$ws.Range("I3").Value = 1564
$ws.Range("J3").Value = 835.5
$ws.Range("K3").Value = 1564
$ws.Range("L3").Value = 835.5
$ws.Range("M3").Value = -1450
$ws.Range("N3").Value = -1063.5

$ws.Range("H105").Value = 1500.5834
$ws.Range("I105").Value = 1414.0952
$ws.Range("K105").Value = 1414.0952
$ws.Range("M105").Value = 332.9048

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2264.01
$ws.Range("I31").Value = 1059.7556
$ws.Range("J31").Value = 3267.5557
$ws.Range("K31").Value = 1059.7556
$ws.Range("L31").Value = 3267.5557
$ws.Range("M31").Value = -764.7556
$ws.Range("N31").Value = -3857.5557

$ws.Range("H34").Value = 2264.01
$ws.Range("I34").Value = 1059.7556
$ws.Range("J34").Value = 3267.5557
$ws.Range("K34").Value = 1059.7556
$ws.Range("L34").Value = 3267.5557
$ws.Range("M34").Value = -857.7556
$ws.Range("N34").Value = -3671.5557

$ws.Range("H134").Value = 1240.7273
$ws.Range("I134").Value = 935
$ws.Range("J134").Value = 3737.5
$ws.Range("K134").Value = 2805
$ws.Range("L134").Value = 11212.5
$ws.Range("M134").Value = -270
$ws.Range("N134").Value = -16282.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 429.25
$ws.Range("I23").Value = 388.8
$ws.Range("J23").Value = 475.92307
$ws.Range("K23").Value = 1166.4
$ws.Range("L23").Value = 1427.76921
$ws.Range("M23").Value = -931.4000000000001
$ws.Range("N23").Value = -1897.76921

$ws.Range("H97").Value = 516.55554
$ws.Range("I97").Value = 260
$ws.Range("J97").Value = 589.8570999999999
$ws.Range("K97").Value = 780
$ws.Range("L97").Value = 1769.5713
$ws.Range("M97").Value = -284
$ws.Range("N97").Value = -2761.5713

$ws.Range("H131").Value = 2865.1526
$ws.Range("I131").Value = 11497.556
$ws.Range("J131").Value = 1311.32
$ws.Range("K131").Value = 34492.66800000001
$ws.Range("L131").Value = 3933.96
$ws.Range("M131").Value = -29452.66800000001
$ws.Range("N131").Value = -14013.96

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4693.4375
$ws.Range("I80").Value = 4759.5
$ws.Range("J80").Value = 4583.3335
$ws.Range("K80").Value = 4759.5
$ws.Range("L80").Value = 4583.3335
$ws.Range("M80").Value = -3761.5
$ws.Range("N80").Value = -6579.3335

$ws.Range("H83").Value = 4693.4375
$ws.Range("I83").Value = 4759.5
$ws.Range("J83").Value = 4583.3335
$ws.Range("K83").Value = 23797.5
$ws.Range("L83").Value = 22916.6675
$ws.Range("M83").Value = -18805.5
$ws.Range("N83").Value = -32900.6675

$ws.Range("H102").Value = 1270.5454
$ws.Range("I102").Value = 1314.3684
$ws.Range("J102").Value = 993
$ws.Range("K102").Value = 1314.3684
$ws.Range("L102").Value = 993
$ws.Range("M102").Value = 307.6315999999999
$ws.Range("N102").Value = -4237

$ws.Range("H113").Value = 1505.0588
$ws.Range("I113").Value = 1791.6666
$ws.Range("J113").Value = 1348.7273
$ws.Range("K113").Value = 1791.6666
$ws.Range("L113").Value = 1348.7273
$ws.Range("M113").Value = 378.3334
$ws.Range("N113").Value = -5688.7273

$ws.Range("H123").Value = 14410.714
$ws.Range("J123").Value = 14410.714
$ws.Range("L123").Value = 14410.714
$ws.Range("N123").Value = -19310.714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3416.6667
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 3900
$ws.Range("K68").Value = 1000
$ws.Range("L68").Value = 3900
$ws.Range("M68").Value = -251
$ws.Range("N68").Value = -5398

$ws.Range("H71").Value = 3416.6667
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 3900
$ws.Range("K71").Value = 5000
$ws.Range("L71").Value = 19500
$ws.Range("M71").Value = -1256
$ws.Range("N71").Value = -26988

$ws.Range("H82").Value = 11906676
$ws.Range("I82").Value = 1849.75
$ws.Range("K82").Value = 1849.75
$ws.Range("M82").Value = -1488.75

$ws.Range("H85").Value = 11906676
$ws.Range("I85").Value = 1849.75
$ws.Range("K85").Value = 1849.75
$ws.Range("M85").Value = -601.75

$ws.Range("H122").Value = 2167.5
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 2502.5
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 7507.5
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -12407.5

$ws.Range("H132").Value = 2012.8309
$ws.Range("I132").Value = 1437.125
$ws.Range("K132").Value = 4311.375
$ws.Range("M132").Value = -1781.375

$ws.Range("H135").Value = 40000
$ws.Range("J135").Value = 40000
$ws.Range("L135").Value = 40000
$ws.Range("N135").Value = -50140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3169.1052
$ws.Range("I62").Value = 2608.25
$ws.Range("J62").Value = 3318.6667
$ws.Range("K62").Value = 2608.25
$ws.Range("L62").Value = 3318.6667
$ws.Range("M62").Value = -1984.25
$ws.Range("N62").Value = -4566.6667

$ws.Range("H65").Value = 3169.1052
$ws.Range("I65").Value = 2608.25
$ws.Range("J65").Value = 3318.6667
$ws.Range("K65").Value = 13041.25
$ws.Range("L65").Value = 16593.3335
$ws.Range("M65").Value = -9921.25
$ws.Range("N65").Value = -22833.3335

$ws.Range("H132").Value = 1513.4
$ws.Range("I132").Value = 1281.2667
$ws.Range("J132").Value = 2209.8
$ws.Range("K132").Value = 3843.800099999999
$ws.Range("L132").Value = 6629.400000000001
$ws.Range("M132").Value = -1313.800099999999
$ws.Range("N132").Value = -11689.4

$ws.Range("H136").Value = 23722.205
$ws.Range("I136").Value = 31833.781
$ws.Range("K136").Value = 95501.34299999999
$ws.Range("M136").Value = -92951.34299999999
